$wb = $excel.ActiveWorkbook

# The same three rows were updated on both the "展览" sheet and the
# "全部类型" sheet (the latter aggregates all entries), so apply the
# identical F-column ("想去人数") updates to each.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 6640
    $ws.Range("F5").Value = 1036
    $ws.Range("F6").Value = 132
}
